# New crime data collected - update the weekly CompStat report
# (Volume/Number header, reporting week dates, and the Crime Complaints table)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text: "Volume 30   Number  5" -> "...Number  6"
# and the reporting-week sentence "... 1/30/2023  Through  2/5/2023" -> "2/6/2023 .. 2/12/2023"
# These are built from several rich-text runs concatenated in one cell each.
# ---------------------------------------------------------------------------
$ws.Range("A8").Characters(21, 1).Text = "6"
$ws.Range("C9").Characters(27, 9).Text = "2/6/2023"
$ws.Range("C9").Characters(47, 8).Text = "2/12/2023"

# ---------------------------------------------------------------------------
# Cells that change TEXT TYPE (numeric -> the literal placeholder text used
# throughout this sheet for zero / undefined-percentage cells). Copying an
# existing placeholder cell (value + format) keeps the exact same shared
# string + style the rest of the sheet uses, instead of minting a new style.
# ---------------------------------------------------------------------------
$ws.Range("D22").Copy($ws.Range("C22"))     # Transit / Murder column -> "0"
$ws.Range("C23").Copy($ws.Range("D18"))     # Burglary D -> "0"
$ws.Range("E23").Copy($ws.Range("E18"))     # Burglary E -> "***.*"

# ---------------------------------------------------------------------------
# Cells that change FROM the placeholder text back TO real numbers. Copy the
# number format from a same-row numeric cell first (so the style matches the
# rest of the numeric columns), then write the real value.
# ---------------------------------------------------------------------------
$ws.Range("F27").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = 2

$ws.Range("G27").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").Value = 1

$ws.Range("H27").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E27").Value = 100

# ---------------------------------------------------------------------------
# Plain numeric value updates (weekly / 28-day / YTD counts and %-change
# figures recomputed for the new reporting week).
# ---------------------------------------------------------------------------
# Row 15 (Rape)
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 2
$ws.Range("K15").Value = 0
$ws.Range("M15").Value = -33.333333333333
$ws.Range("N15").Value = -60

# Row 16 (Robbery)
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -25
$ws.Range("F16").Value = 14
$ws.Range("H16").Value = 16.666666666666
$ws.Range("I16").Value = 21
$ws.Range("J16").Value = 17
$ws.Range("K16").Value = 23.529411764705
$ws.Range("L16").Value = 31.25
$ws.Range("M16").Value = -46.153846153846
$ws.Range("N16").Value = -82.786885245901

# Row 17 (Fel. Assault)
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -20
$ws.Range("F17").Value = 19
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = -5
$ws.Range("I17").Value = 36
$ws.Range("J17").Value = 29
$ws.Range("K17").Value = 24.137931034482
$ws.Range("L17").Value = 38.461538461538
$ws.Range("M17").Value = 157.142857142857
$ws.Range("N17").Value = -7.692307692307

# Row 18 (Burglary)
$ws.Range("C18").Value = 1
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 57.142857142857
$ws.Range("I18").Value = 17
$ws.Range("K18").Value = 54.545454545454
$ws.Range("L18").Value = 41.666666666666
$ws.Range("M18").Value = -37.037037037037
$ws.Range("N18").Value = -89.570552147239

# Row 19 (Gr. Larceny)
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 23
$ws.Range("E19").Value = -65.217391304347
$ws.Range("F19").Value = 51
$ws.Range("G19").Value = 71
$ws.Range("H19").Value = -28.169014084507
$ws.Range("I19").Value = 74
$ws.Range("J19").Value = 104
$ws.Range("K19").Value = -28.846153846153
$ws.Range("L19").Value = 48
$ws.Range("M19").Value = 111.428571428571
$ws.Range("N19").Value = 8.823529411764

# Row 20 (G.L.A.)
$ws.Range("C20").Value = 10
$ws.Range("D20").Value = 11
$ws.Range("E20").Value = -9.090909090909
$ws.Range("G20").Value = 31
$ws.Range("H20").Value = -16.129032258064
$ws.Range("I20").Value = 42
$ws.Range("J20").Value = 40
$ws.Range("K20").Value = 5
$ws.Range("L20").Value = 133.333333333333
$ws.Range("M20").Value = 10.526315789473
$ws.Range("N20").Value = -90.344827586206

# Row 21 (TOTAL)
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 43
$ws.Range("E21").Value = -39.534883720930
$ws.Range("F21").Value = 122
$ws.Range("G21").Value = 142
$ws.Range("H21").Value = -14.084507042253
$ws.Range("I21").Value = 192
$ws.Range("J21").Value = 203
$ws.Range("K21").Value = -5.418719211822
$ws.Range("L21").Value = 57.377049180327
$ws.Range("M21").Value = 23.076923076923
$ws.Range("N21").Value = -77.060931899641

# Row 22 (Transit)
$ws.Range("M22").Value = 33.333333333333

# Row 24 (Petit Larceny)
$ws.Range("C24").Value = 24
$ws.Range("D24").Value = 33
$ws.Range("E24").Value = -27.272727272727
$ws.Range("F24").Value = 92
$ws.Range("G24").Value = 93
$ws.Range("H24").Value = -1.075268817204
$ws.Range("I24").Value = 135
$ws.Range("J24").Value = 129
$ws.Range("K24").Value = 4.651162790697
$ws.Range("L24").Value = 64.634146341463
$ws.Range("M24").Value = 90.140845070422

# Row 25 (Misd. Assault)
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 66.666666666666
$ws.Range("F25").Value = 38
$ws.Range("G25").Value = 28
$ws.Range("H25").Value = 35.714285714285
$ws.Range("I25").Value = 59
$ws.Range("J25").Value = 41
$ws.Range("K25").Value = 43.902439024390
$ws.Range("L25").Value = 34.090909090909
$ws.Range("M25").Value = 22.916666666666

# Row 26 (UCR Rape*)
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 3
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 50

# Row 27 (Other Sex Crimes) - remaining numeric updates
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 12
$ws.Range("J27").Value = 6
$ws.Range("L27").Value = 50

# Row 28 (Shooting Vic.)
$ws.Range("N28").Value = -71.428571428571

# Row 29 (Shooting Inc.)
$ws.Range("N29").Value = -85.714285714285
